{"js": "// Update the date heading and the 25 division-equation table cells to match\n// the new day's worksheet content.\n\n// 1) Update the date paragraph (first paragraph in the body, before the table).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-06-30 Sunday\", \"Replace\");\n\n// 2) Update the equation cells inside the single table. The table has 20\n// rows x 5 columns; only rows 0, 4, 8, 12, 16 hold equation text, the rest\n// are blank spacer rows. Replacements are applied by (row, col) position\n// (not by text-match) because several source equations are duplicated\n// (e.g. \"49\u00f72=\" appears twice) and some new values collide with other\n// original values (e.g. \"97\u00f79=\" / \"32\u00f72=\"), so a sequential find/replace\n// would corrupt later cells.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = [\n  [0, 0, \"41\u00f76=\"],\n  [0, 1, \"33\u00f78=\"],\n  [0, 2, \"78\u00f78=\"],\n  [0, 3, \"77\u00f77=\"],\n  [0, 4, \"56\u00f77=\"],\n  [4, 0, \"12\u00f72=\"],\n  [4, 1, \"22\u00f75=\"],\n  [4, 2, \"15\u00f73=\"],\n  [4, 3, \"97\u00f79=\"],\n  [4, 4, \"70\u00f72=\"],\n  [8, 0, \"72\u00f72=\"],\n  [8, 1, \"33\u00f72=\"],\n  [8, 2, \"10\u00f79=\"],\n  [8, 3, \"32\u00f72=\"],\n  [8, 4, \"17\u00f76=\"],\n  [12, 0, \"31\u00f76=\"],\n  [12, 1, \"93\u00f73=\"],\n  [12, 2, \"57\u00f78=\"],\n  [12, 3, \"57\u00f79=\"],\n  [12, 4, \"40\u00f74=\"],\n  [16, 0, \"63\u00f75=\"],\n  [16, 1, \"86\u00f74=\"],\n  [16, 2, \"24\u00f79=\"],\n  [16, 3, \"69\u00f72=\"],\n  [16, 4, \"53\u00f72=\"],\n];\n\nfor (const [row, col, text] of newValues) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division-equation table cells to match\n# the new day's worksheet content.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph in the body, before the table).\n$d.Paragraphs.Item(1).Range.Text = \"2024-06-30 Sunday\"\n\n# 2) Update the equation cells inside the single table. The table has 20\n# rows x 5 columns (1-based in the COM object model); only rows 1, 5, 9, 13,\n# 17 hold equation text, the rest are blank spacer rows. Replacements are\n# applied by (row, col) position (not by text-match) because several source\n# equations are duplicated (e.g. \"49\u00f72=\" appears twice) and some new values\n# collide with other original values (e.g. \"97\u00f79=\" / \"32\u00f72=\"), so a\n# sequential find/replace would corrupt later cells.\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(1, 1, \"41\u00f76=\"),\n    @(1, 2, \"33\u00f78=\"),\n    @(1, 3, \"78\u00f78=\"),\n    @(1, 4, \"77\u00f77=\"),\n    @(1, 5, \"56\u00f77=\"),\n    @(5, 1, \"12\u00f72=\"),\n    @(5, 2, \"22\u00f75=\"),\n    @(5, 3, \"15\u00f73=\"),\n    @(5, 4, \"97\u00f79=\"),\n    @(5, 5, \"70\u00f72=\"),\n    @(9, 1, \"72\u00f72=\"),\n    @(9, 2, \"33\u00f72=\"),\n    @(9, 3, \"10\u00f79=\"),\n    @(9, 4, \"32\u00f72=\"),\n    @(9, 5, \"17\u00f76=\"),\n    @(13, 1, \"31\u00f76=\"),\n    @(13, 2, \"93\u00f73=\"),\n    @(13, 3, \"57\u00f78=\"),\n    @(13, 4, \"57\u00f79=\"),\n    @(13, 5, \"40\u00f74=\"),\n    @(17, 1, \"63\u00f75=\"),\n    @(17, 2, \"86\u00f74=\"),\n    @(17, 3, \"24\u00f79=\"),\n    @(17, 4, \"69\u00f72=\"),\n    @(17, 5, \"53\u00f72=\")\n)\n\nforeach ($entry in $newValues) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $text = $entry[2]\n    $t.Cell($row, $col).Range.Text = $text\n}\n"}
